# Update on 2018-06-29, 支出生活费400
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 was a blank placeholder row (styled but empty). Populate it with a
# new expense entry: #41, 支出 (expense), 400, 2018-06-19, 生活费 (category),
# remark "生活费(6/20-6/30)".
# Copy the formatting from the row above (row 42) first so the new row picks
# up the exact same cell styles (fill/border/number-format) as the rest of
# the table, then overwrite the values.
$ws.Range("B42:G42").Copy($ws.Range("B43:G43")) | Out-Null

$ws.Range("B43").Value = 41
$ws.Range("C43").Value = "支出"
$ws.Range("D43").Value = 400
$ws.Range("E43").Value = "6/19/2018"
$ws.Range("F43").Value = "生活费"
$ws.Range("G43").Value = "生活费(6/20-6/30)"

# Update the view selection to match the edited area.
$ws.Range("N51").Select() | Out-Null
